$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("llm_extension_with_qa_extension")

# Big multi-line prompt text for the new row (single-quoted here-string -> no escaping)
$promptText = @'
`
    ## Task: Generate knowledge graph query commands for Sparklis.
    ## Format:  
    1. Think step by step about what entities and relationships are needed 
    2. Then finish your response by a list of commands, separated by semicolons (;), and wrapped in <commands>...</commands>.  
    ### Available Commands:
    - a [concept] → Retrieve entities of a given concept (e.g., "a book" to find books).
    - [entity] → Retrieve an entity (e.g., "Albert Einstein" to find the entity representing Einstein).
    - forwardProperty [property] → Filter by property (e.g., "forwardProperty director" to find films directed by someone).
    - backwardProperty [property] → Reverse relation (e.g., "backwardProperty director" to find directors of films).
    - higherThan [number], lowerThan [number] → Value constraints.
    - after [date], before [date] → Time constraints.  
    - and, or → Logical operators.  
    ## Examples:
    Q: At which school went Yayoi Kusama?
    A: Starting from the list of entities named Yayoi Kusama seems the best approach. Then, I just need to find the relationship that represents at which school she was educated.
    <commands>Yayoi Kusama ; forwardProperty educated at</commands> 
    Q: What is the boiling point of water?
    A: The core of the request is WATER. From this entity I will probably be able to get a property such as its BOILING POINT.  
    <commands>water; forwardProperty boiling point</commands>  
    Q: Movies by Spielberg or Tim Burton after 1980?
    A: I need to find FILMS by Spielberg or Burton released after 1980. I can start by listing FILMS and then filter by DIRECTOR and RELEASE DATE. 
    <commands>a film; forwardProperty director; Tim Burton; or; Spielberg; forwardProperty release date; after 1980</commands>  
    Q: among the founders of tencent company, who has been member of national people' congress?"
    A: I can start by finding FOUNDERS of something called TENCENT. Then, I can filter by people who have been members of the NATIONAL PEOPLE'S CONGRESS.
    <commands>backwardProperty founder of; Tencent ; forwardProperty position ; National People's Congress</commands>
    `
'@

# Add the new data row (row 10). Cells are written C, B, A in that order so that the
# shared-string table gets new entries in the same order as the target file
# (47 = prompt text, 48 = "patch backward", 49 = "v10").
$ws1.Range("C10").Value = $promptText
$ws1.Range("B10").Value = "patch backward"
$ws1.Range("A10").Value = "v10"

# Column C uses the workbook's wrap-text cell style (same as the rest of column C).
$ws1.Range("C10").WrapText = $true
$ws1.Rows.Item(10).RowHeight = 409.5

# This sheet becomes the active / selected tab, scrolled near the new row, with A11 selected.
$ws1.Activate()
$ws1.Range("A11").Select()

